# Daily attendance processing - rotate "Recorded By" (column G) entries.
# For each data row, the comma-separated list of recorders in column G is
# right-rotated by one position (the last entry moves to the front).
# Lists with 0 or 1 entries are left unchanged (rotation is a no-op for them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -gt 1) {
        $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
        $cell.Value = [string]::Join(", ", $rotated)
    }
}
